$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "600px" breakpoint row: drop the highlighted ("Neutral") style so it
# matches the plain look used by the rest of the table (E5, E6, E7).
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").HorizontalAlignment = -4131

# Row 3 ("tablet landscape"): same style fix, and replace the numeric 304 value
# with the new annotated text "320 (228)".
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").HorizontalAlignment = -4131
$ws.Range("E3").Value = "320 (228)"

# Row 4 ("laptop"): same style fix, and correct the numeric value from 304 to 320.
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").HorizontalAlignment = -4131
$ws.Range("E4").Value = 320

# Update the saved cursor/selection position on the sheet.
$ws.Range("E8").Select()
